# Update need_to_buy.xlsx per the latest R run (columns B, C, F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7: TISG_PDR_G (B), fcs (C) and need_to_buy_MW (F) are now blank
# (the refresh from R produced NA there). Clear the values but keep the
# cell entries present in the sheet (touching Style keeps the <c> node
# alive instead of Excel pruning the now-empty cell on save).
$ws.Range("B2:C7").ClearContents()
$ws.Range("B2:C7").Style = "Normal"
$ws.Range("F2:F7").ClearContents()
$ws.Range("F2:F7").Style = "Normal"

# Rows 8-15: refreshed B/C/F values from the new R run.
$updates = @(
    @{ Row = 8;  B = 5230.43751372964; C = 4762.76196774228; F = 137.10596720886 },
    @{ Row = 9;  B = 5230.43751372964; C = 5010.51019936627; F = 147.428810193193 },
    @{ Row = 10; B = 5230.43751372964; C = 5212.39493853064; F = 155.840674325042 },
    @{ Row = 11; B = 5230.43751372964; C = 5812.90607844431; F = 180.861971821445 },
    @{ Row = 12; B = 5230.43751372964; C = 5337.7721050942;  F = 161.064722931857 },
    @{ Row = 13; B = 1198.04912219854; C = 3211.54743728158; F = 57.62422775346 },
    @{ Row = 14; B = 1070.29084929054; C = 3178.52433118836; F = 55.8169678707426 },
    @{ Row = 15; B = 5457.44872542722; C = 5928.99508284696; F = 189.632203142489 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
